# Final update of flow charts to be vector based.
# Update of inter-vocab relationship section
#
# The whole flow-chart diagram on slide 4 (all top-level shapes) is shifted
# by a constant offset: (-1296177, -1615207) EMU.
# PowerPoint COM exposes Shape.Left / Shape.Top in points (1 pt = 12700 EMU),
# so we convert the EMU delta to points with full precision and apply it to
# every shape on the slide (including groups, connectors and pictures).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$dxEmu = -1296177
$dyEmu = -1615207
$dxPt = $dxEmu / 12700.0
$dyPt = $dyEmu / 12700.0

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $sh.Left = $sh.Left + $dxPt
    $sh.Top = $sh.Top + $dyPt
}
